# Apply updated dSF (column F) values for the stout_eric 2022 workbook.
# Column F holds "dSF" per the header row; rows 2-25 correspond to
# A column index 0-23 (one row per data point).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -12
    3  = -4
    4  = -5
    5  = -4
    6  = 1
    9  = -7
    10 = -3
    11 = 0
    14 = -1
    16 = -3
    21 = -7
    25 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
